$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.678.56"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.106.97"
$ws.Range("E3").Value = "  +10.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5204"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4399"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08962"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("D13").Value = "2.106.79"
$ws.Range("E13").Value = "  +10.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.791"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.654"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001134"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06607"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.411"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.09%  "
$ws.Range("D23").Value = "30.833.41"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("D25").Value = "2.360.03"
$ws.Range("E25").Value = "  +10.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.260"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.532"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.203"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.915"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.525"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +27.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02589"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.619"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06764"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.518"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2244"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6773"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.252"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.04%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6304"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.251"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.648"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.273"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.90%  "
